# EN325_WFS_Panels.xlsx edit
# - rename Sheet2 -> exportPanels
# - add new sheet exportPolygon (pulls a polygon point list out of Sheet1)
# - change the W-offset default (J31) from 0.4 to 0 on Sheet1 (cascades through
#   the O/P/Q/R "NORM" columns + Sheet2 export columns via existing formulas)
# - add a small "Polygon" / "wOffset" table (rows 35-44) on Sheet1 that feeds
#   the new exportPolygon sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheets: rename Sheet2, add exportPolygon right after it
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "exportPanels"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "exportPolygon"

# ---------------------------------------------------------------------------
# 2. Sheet1: core value change driving the whole recalculation
# ---------------------------------------------------------------------------
$ws1.Range("J31").Value = 0

# ---------------------------------------------------------------------------
# 3. Sheet1: new "Polygon" / point table, rows 35-44
# ---------------------------------------------------------------------------
$ws1.Range("A35:C35").Merge()
$ws1.Range("A35").Value = "Polygon"

$ws1.Range("E35:F35").Merge()
$ws1.Range("E35").Value = "wOffset"

$ws1.Range("A36").Value = "Point"
$ws1.Range("B36").Value = "x1"
$ws1.Range("C36").Value = "y1"
$ws1.Range("E36").Value = "x1"
$ws1.Range("F36").Value = "y1"

# formatting: reuse existing styles instead of inventing new ones
$ws1.Range("B1").Copy()
$ws1.Range("E35:F35").PasteSpecial(-4122)

$ws1.Range("J31").Copy()
$ws1.Range("A36:C36").PasteSpecial(-4122)
$ws1.Range("E36:F36").PasteSpecial(-4122)

$ws1.Range("B1").Copy()
$ws1.Range("A35:C35").PasteSpecial(-4122)
$ws1.Range("A35:C35").Font.Bold = $true

$ws1.Range("A37").Value = 1
$ws1.Range("B37").Formula = "=B25"
$ws1.Range("C37").Formula = "=C25"
$ws1.Range("E37").Formula = "=B37+`$J`$31"
$ws1.Range("F37").Formula = "=C37 + `$K`$31"

$ws1.Range("A38").Value = 2
$ws1.Range("B38").Formula = "=D4"
$ws1.Range("C38").Formula = "=E4"
$ws1.Range("E38").Formula = "=B38+`$J`$31"
$ws1.Range("F38").Formula = "=C38 + `$K`$31"

$ws1.Range("A39").Value = 3
$ws1.Range("B39").Formula = "=D6"
$ws1.Range("C39").Formula = "=E6"
$ws1.Range("E39").Formula = "=B39+`$J`$31"
$ws1.Range("F39").Formula = "=C39 + `$K`$31"

$ws1.Range("A40").Value = 4
$ws1.Range("B40").Formula = "=D10"
$ws1.Range("C40").Formula = "=E10"
$ws1.Range("E40").Formula = "=B40+`$J`$31"
$ws1.Range("F40").Formula = "=C40 + `$K`$31"

$ws1.Range("A41").Value = 5
$ws1.Range("B41").Formula = "=D12"
$ws1.Range("C41").Formula = "=E12"
$ws1.Range("E41").Formula = "=B41+`$J`$31"
$ws1.Range("F41").Formula = "=C41 + `$K`$31"

$ws1.Range("A42").Value = 6
$ws1.Range("B42").Formula = "=D16"
$ws1.Range("C42").Formula = "=E16"
$ws1.Range("E42").Formula = "=B42+`$J`$31"
$ws1.Range("F42").Formula = "=C42 + `$K`$31"

$ws1.Range("A43").Value = 7
$ws1.Range("B43").Formula = "=D18"
$ws1.Range("C43").Formula = "=E18"
$ws1.Range("E43").Formula = "=B43+`$J`$31"
$ws1.Range("F43").Formula = "=C43 + `$K`$31"

$ws1.Range("A44").Value = 8

# ---------------------------------------------------------------------------
# 4. exportPolygon sheet: pull the point table back out via formulas
# ---------------------------------------------------------------------------
$ws3.Range("A1").Formula = "=Sheet1!A36"
$ws3.Range("B1").Formula = "=Sheet1!B36"
$ws3.Range("C1").Formula = "=Sheet1!C36"

$rows = 2..8
$srcRow = 37
foreach ($r in $rows) {
    $ws3.Range("A$r").Formula = "=Sheet1!A$srcRow"
    $ws3.Range("B$r").Formula = "=Sheet1!E$srcRow"
    $ws3.Range("C$r").Formula = "=Sheet1!F$srcRow"
    $srcRow++
}

$ws1.Range("V3").Copy()
$ws3.Range("B2:C8").PasteSpecial(-4122)

# restore formula content after the format-only paste overwrote it
$srcRow = 37
foreach ($r in $rows) {
    $ws3.Range("B$r").Formula = "=Sheet1!E$srcRow"
    $ws3.Range("C$r").Formula = "=Sheet1!F$srcRow"
    $srcRow++
}

$ws3.Range("A1").Select()

# ---------------------------------------------------------------------------
# 5. Cosmetic: selection on Sheet1
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("J32").Select()

$excel.CutCopyMode = $false
